$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "Status" result cell (D2) - the template no longer emits this value.
$ws.Range("D2").Value = $null

# Move the active selection from C2 to D2, matching the updated template.
$ws.Range("D2").Select()
